$d = $word.ActiveDocument

# Locate the paragraph text "iaulne &" and split off the leading "i",
# turning it into "j" in its own run (without the color formatting)
# while the remainder "aulne &" keeps the original color formatting.
$rng = $d.Content
[void]$rng.Find.Execute("iaulne &", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$start = $rng.Start
$charRange = $d.Range($start, $start + 1)
[void]$charRange.Delete()

$insPoint = $d.Range($start, $start)
$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">j</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$insPoint.InsertXML($xmlFrag)
